# Add cantrals by cantons: fill in the previously-zeroed "Bâle-Ville" (column M)
# counts for years 2006-2017 (rows 2-13) on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 2
$ws.Range("M3").Value = 3
$ws.Range("M4").Value = 3
$ws.Range("M5").Value = 3
$ws.Range("M6").Value = 3
$ws.Range("M7").Value = 3
$ws.Range("M8").Value = 3
$ws.Range("M9").Value = 3
$ws.Range("M10").Value = 3
$ws.Range("M11").Value = 3
$ws.Range("M12").Value = 3
$ws.Range("M13").Value = 4

# Match the author's final selection in the saved workbook.
$ws.Range("M3").Select()
